$wb = $excel.ActiveWorkbook

# --- "region" sheet: rename the "JAVA" region-group label to "JAWA" ---
$regionSheet = $wb.Worksheets.Item("region")
for ($r = 1; $r -le 36; $r++) {
    $cell = $regionSheet.Cells.Item($r, 2)
    $v = $cell.Value()
    if ($v -eq "JAVA") {
        $cell.Value = "JAWA"
    }
}

# --- Switch the active sheet/selection back to "region" ---
$regionSheet.Activate()
$regionSheet.Range("B25").Select()
